$d = $word.ActiveDocument

$d.Content.Find.Execute("108×8=864", $true, $false, $false, $false, $false, $true, 1, $false, "262×5=1310", 2) | Out-Null
$d.Content.Find.Execute("516×6=3096", $true, $false, $false, $false, $false, $true, 1, $false, "226×3=678", 2) | Out-Null
$d.Content.Find.Execute("897×4=3588", $true, $false, $false, $false, $false, $true, 1, $false, "111×6=666", 2) | Out-Null
$d.Content.Find.Execute("528×6=3168", $true, $false, $false, $false, $false, $true, 1, $false, "187×6=1122", 2) | Out-Null
$d.Content.Find.Execute("919×4=3676", $true, $false, $false, $false, $false, $true, 1, $false, "903×6=5418", 2) | Out-Null
$d.Content.Find.Execute("858×2=1716", $true, $false, $false, $false, $false, $true, 1, $false, "300×3=900", 2) | Out-Null
$d.Content.Find.Execute("534×9=4806", $true, $false, $false, $false, $false, $true, 1, $false, "789×6=4734", 2) | Out-Null
$d.Content.Find.Execute("692×3=2076", $true, $false, $false, $false, $false, $true, 1, $false, "896×2=1792", 2) | Out-Null
$d.Content.Find.Execute("993×4=3972", $true, $false, $false, $false, $false, $true, 1, $false, "990×5=4950", 2) | Out-Null
$d.Content.Find.Execute("885×4=3540", $true, $false, $false, $false, $false, $true, 1, $false, "934×8=7472", 2) | Out-Null
$d.Content.Find.Execute("779×7=5453", $true, $false, $false, $false, $false, $true, 1, $false, "676×3=2028", 2) | Out-Null
$d.Content.Find.Execute("725×3=2175", $true, $false, $false, $false, $false, $true, 1, $false, "389×3=1167", 2) | Out-Null
$d.Content.Find.Execute("688×3=2064", $true, $false, $false, $false, $false, $true, 1, $false, "992×7=6944", 2) | Out-Null
$d.Content.Find.Execute("279×5=1395", $true, $false, $false, $false, $false, $true, 1, $false, "584×3=1752", 2) | Out-Null
$d.Content.Find.Execute("927×9=8343", $true, $false, $false, $false, $false, $true, 1, $false, "194×5=970", 2) | Out-Null
$d.Content.Find.Execute("564×3=1692", $true, $false, $false, $false, $false, $true, 1, $false, "647×9=5823", 2) | Out-Null
$d.Content.Find.Execute("800×5=4000", $true, $false, $false, $false, $false, $true, 1, $false, "955×2=1910", 2) | Out-Null
$d.Content.Find.Execute("543×6=3258", $true, $false, $false, $false, $false, $true, 1, $false, "489×8=3912", 2) | Out-Null
$d.Content.Find.Execute("101×6=606", $true, $false, $false, $false, $false, $true, 1, $false, "482×5=2410", 2) | Out-Null
$d.Content.Find.Execute("979×5=4895", $true, $false, $false, $false, $false, $true, 1, $false, "528×4=2112", 2) | Out-Null
$d.Content.Find.Execute("385×5=1925", $true, $false, $false, $false, $false, $true, 1, $false, "483×2=966", 2) | Out-Null
$d.Content.Find.Execute("426×9=3834", $true, $false, $false, $false, $false, $true, 1, $false, "258×4=1032", 2) | Out-Null
$d.Content.Find.Execute("198×5=990", $true, $false, $false, $false, $false, $true, 1, $false, "600×4=2400", 2) | Out-Null
$d.Content.Find.Execute("284×4=1136", $true, $false, $false, $false, $false, $true, 1, $false, "822×6=4932", 2) | Out-Null
$d.Content.Find.Execute("558×7=3906", $true, $false, $false, $false, $false, $true, 1, $false, "128×2=256", 2) | Out-Null

$d.Save()
